$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (23-nov) before column DX ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("DX").Insert()
$ws1.Range("DX1").Value = "23-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 128).Value = "-"
}

# --- Sheet "Gaz": append new row for 2025-11-21 ---
# (NumberFormat "@" forces the date-looking string to be kept as literal
#  text instead of being auto-parsed into a date serial; ClearFormats()
#  afterwards drops the cell back to the default/unstyled state so the
#  cell ends up as a plain text cell, matching the rest of column A.)
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A157").NumberFormat = "@"
$ws2.Range("A157").Value = "2025-11-21"
$ws2.Range("A157").ClearFormats()
$ws2.Range("B157").Value = 29.35

# --- Sheet "CO2": append new row for 2025-11-21 ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A157").NumberFormat = "@"
$ws3.Range("A157").Value = "2025-11-21"
$ws3.Range("A157").ClearFormats()
$ws3.Range("B157").Value = 80.28
